$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Investment_Cost": fill in Diesel_storage (row 22) and add a new row
# for Ammonia_storage (row 26), plus two blank table rows (27-28) and a
# totals row (29) via the table's ShowTotals feature.
# ---------------------------------------------------------------------------
$wsCost = $wb.Worksheets.Item("Investment_Cost")

# Row 22 - Diesel_storage: fill in the investment cost values.
$wsCost.Range("B22:F22").Value = 25
$wsCost.Range("H22").Value = 0.05
$wsCost.Range("I22").Value = "Average value taken from chart and calculated from `$/L to €/MWh using the energy density"

# Row 26 - new Ammonia_storage entry (same number formats as the rest of
# the table body: #,##0.00 for the cost/notes columns, 0% for the share).
$wsCost.Range("A26").Value = "Ammonia_storage"
$wsCost.Range("B26:F26").NumberFormat = "#,##0.00"
$wsCost.Range("B26:F26").Value = 4400
$wsCost.Range("H26").NumberFormat = "0%"
$wsCost.Range("I26").NumberFormat = "#,##0.00"
$wsCost.Range("I26").Value = "Costs are for 2010 but book is from 2024 so these are probably still relevant, Tanks are either pressurized or refrigerated so power supply needed (https://www.sciencedirect.com/science/article/pii/B9780323885164000111)"

# Rows 27-28 stay blank but still carry the table body's number formats.
foreach ($r in 27, 28) {
    $wsCost.Range("B$r`:F$r").NumberFormat = "#,##0.00"
    $wsCost.Range("H$r").NumberFormat = "0%"
    $wsCost.Range("I$r").NumberFormat = "#,##0.00"
}

# Resize the table to cover the new body rows, then turn on the totals row
# (this is what pushes the table ref from A1:I25 to A1:I29 with
# totalsRowCount="1" / autoFilter A1:I28, matching a normal Excel "add rows
# then show totals" edit). Pre-format row 29 the same way so the totals row
# cells carry the expected number formats.
$wsCost.Range("B29:F29").NumberFormat = "#,##0.00"
$wsCost.Range("H29").NumberFormat = "0%"
$wsCost.Range("I29").NumberFormat = "#,##0.00"

$lo = $wsCost.ListObjects.Item("Table1")
$lo.Resize($wsCost.Range("A1:I28")) | Out-Null
$lo.ShowTotals = $true

# Restore view state roughly matching the target (scrolled down, C26
# selected).
$wsCost.Activate()
try { $excel.ActiveWindow.ScrollRow = 12 } catch {}
$wsCost.Range("C26").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Sources": add the two new source notes with hyperlinks.
# ---------------------------------------------------------------------------
$wsSrc = $wb.Worksheets.Item("Sources")

$wsSrc.Range("A6").Value = "Ammonia storage: "
$wsSrc.Range("C6").Value = "https://www.sciencedirect.com/science/article/pii/B9780323885164000111, Figure 11.4"
$wsSrc.Hyperlinks.Add($wsSrc.Range("C6"), "https://www.sciencedirect.com/science/article/pii/B9780323885164000111") | Out-Null

$wsSrc.Range("A7").Value = "Diesel storage:"
$wsSrc.Range("C7").Value = "https://thundersaidenergy.com/downloads/storage-tank-costs-storing-oil-energy-water-and-chemicals/"
$wsSrc.Hyperlinks.Add($wsSrc.Range("C7"), "https://thundersaidenergy.com/downloads/storage-tank-costs-storing-oil-energy-water-and-chemicals/") | Out-Null

$wsSrc.Activate()
$wsSrc.Range("C8").Select() | Out-Null

$wsCost.Activate()
